# Weekly update: a new price record (week of 2023-11-06) is inserted
# above the existing data, pushing the previous rows 57-94 down to 58-95.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 57 (shifts old rows 57:94 down to 58:95)
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record
$ws.Cells.Item(57, 1).Value = 7
$ws.Cells.Item(57, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(57, 3).Value = "Ñuble"
$ws.Cells.Item(57, 4).Value = 45236
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(57, 6).Value = 100112026
$ws.Cells.Item(57, 7).Value = "Haba"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 60
$ws.Cells.Item(57, 11).Value = 10000
$ws.Cells.Item(57, 12).Value = 10000
$ws.Cells.Item(57, 13).Value = 10000
$ws.Cells.Item(57, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(57, 16).Value = 400
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
